# Updated cryptos list on Fri Sep 15 17:39:25 UTC 2023 with GitHub Actions
#
# Refreshes the Price / Volume(1h) columns for the existing coin rows and
# applies the row 48-51 re-ranking: a new "BabyDogeCoin" entry lands at row
# 48, shifting Cronos / EnergySwap / Mantle down one row each; USDD
# (previously row 51) drops off the bottom of the A1:E51 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> 1-based index (A=1 .. E=5)
$colIndex = @{ 'A' = 1; 'B' = 2; 'C' = 3; 'D' = 4; 'E' = 5 }

# Each entry: row R, column letter C, new value V. Text=$true forces the
# cell to keep storing V as literal text (instead of Excel auto-converting
# a numeric-looking string like "18.79" into a real number), matching the
# source feed's plain-text price column.
$updates = @(
    @{R=2; C='D'; V='26.442.03'; Text=$false}
    @{R=2; C='E'; V='  -0.85%  '; Text=$false}
    @{R=3; C='D'; V='1.625.77'; Text=$false}
    @{R=3; C='E'; V='  -0.75%  '; Text=$false}
    @{R=4; C='E'; V='  +0.28%  '; Text=$false}
    @{R=5; C='D'; V='212.97'; Text=$true}
    @{R=5; C='E'; V='  -0.04%  '; Text=$false}
    @{R=6; C='D'; V='0.498'; Text=$true}
    @{R=6; C='E'; V='  +1.26%  '; Text=$false}
    @{R=7; C='E'; V='  +0.25%  '; Text=$false}
    @{R=8; C='E'; V='  +0.06%  '; Text=$false}
    @{R=9; C='E'; V='  -1.45%  '; Text=$false}
    @{R=10; C='D'; V='18.79'; Text=$true}
    @{R=10; C='E'; V='  -1.61%  '; Text=$false}
    @{R=11; C='E'; V='  +0.96%  '; Text=$false}
    @{R=12; C='D'; V='1.853.50'; Text=$false}
    @{R=12; C='E'; V='  -0.61%  '; Text=$false}
    @{R=13; C='D'; V='1.638.10'; Text=$false}
    @{R=13; C='E'; V='  -0.03%  '; Text=$false}
    @{R=14; C='D'; V='4.12'; Text=$true}
    @{R=14; C='E'; V='  +1.67%  '; Text=$false}
    @{R=15; C='D'; V='0.521'; Text=$true}
    @{R=15; C='E'; V='  -0.53%  '; Text=$false}
    @{R=16; C='D'; V='64.77'; Text=$true}
    @{R=16; C='E'; V='  +2.72%  '; Text=$false}
    @{R=17; C='D'; V='26.494.07'; Text=$false}
    @{R=17; C='E'; V='  -0.61%  '; Text=$false}
    @{R=18; C='E'; V='  -0.20%  '; Text=$false}
    @{R=19; C='D'; V='214.92'; Text=$true}
    @{R=19; C='E'; V='  +2.20%  '; Text=$false}
    @{R=20; C='E'; V='  +0.23%  '; Text=$false}
    @{R=21; C='E'; V='  -0.45%  '; Text=$false}
    @{R=22; C='D'; V='6.25'; Text=$true}
    @{R=22; C='E'; V='  +1.87%  '; Text=$false}
    @{R=23; C='D'; V='9.27'; Text=$true}
    @{R=23; C='E'; V='  -1.27%  '; Text=$false}
    @{R=24; C='D'; V='2.00'; Text=$true}
    @{R=24; C='E'; V='  +4.20%  '; Text=$false}
    @{R=25; C='D'; V='148.58'; Text=$true}
    @{R=25; C='E'; V='  +1.78%  '; Text=$false}
    @{R=26; C='D'; V='1.01'; Text=$true}
    @{R=26; C='E'; V='  +0.34%  '; Text=$false}
    @{R=27; C='D'; V='0.119'; Text=$true}
    @{R=27; C='E'; V='  -0.26%  '; Text=$false}
    @{R=28; C='D'; V='6.82'; Text=$true}
    @{R=28; C='E'; V='  +1.62%  '; Text=$false}
    @{R=29; C='D'; V='15.54'; Text=$true}
    @{R=29; C='E'; V='  +0.81%  '; Text=$false}
    @{R=30; C='E'; V='  -1.45%  '; Text=$false}
    @{R=31; C='E'; V='  -0.85%  '; Text=$false}
    @{R=32; C='E'; V='  +2.69%  '; Text=$false}
    @{R=33; C='D'; V='2.94'; Text=$true}
    @{R=33; C='E'; V='  -0.66%  '; Text=$false}
    @{R=34; C='E'; V='  -0.92%  '; Text=$false}
    @{R=35; C='D'; V='2.39'; Text=$true}
    @{R=35; C='E'; V='  -0.76%  '; Text=$false}
    @{R=36; C='D'; V='1.219.93'; Text=$false}
    @{R=36; C='E'; V='  +4.38%  '; Text=$false}
    @{R=37; C='E'; V='  +3.93%  '; Text=$false}
    @{R=38; C='E'; V='  +0.22%  '; Text=$false}
    @{R=39; C='D'; V='0.793'; Text=$true}
    @{R=39; C='E'; V='  -2.15%  '; Text=$false}
    @{R=40; C='E'; V='  +0.69%  '; Text=$false}
    @{R=41; C='D'; V='2.26'; Text=$true}
    @{R=41; C='E'; V='  -2.34%  '; Text=$false}
    @{R=42; C='D'; V='0.792'; Text=$true}
    @{R=42; C='E'; V='  -0.42%  '; Text=$false}
    @{R=43; C='E'; V='  -0.10%  '; Text=$false}
    @{R=44; C='D'; V='1.764.58'; Text=$false}
    @{R=44; C='E'; V='  -0.43%  '; Text=$false}
    @{R=45; C='D'; V='93.07'; Text=$true}
    @{R=45; C='E'; V='  +0.83%  '; Text=$false}
    @{R=46; C='E'; V='  +1.35%  '; Text=$false}
    @{R=47; C='D'; V='54.84'; Text=$true}
    @{R=47; C='E'; V='  +0.19%  '; Text=$false}
    @{R=48; C='B'; V='BabyDogeCoin'; Text=$false}
    @{R=48; C='C'; V='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; Text=$false}
    @{R=48; C='D'; V='0.0₆0104'; Text=$false}
    @{R=48; C='E'; V='  -0.79%  '; Text=$false}
    @{R=49; C='B'; V='Cronos'; Text=$false}
    @{R=49; C='C'; V='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Text=$false}
    @{R=49; C='D'; V='0.0509'; Text=$true}
    @{R=49; C='E'; V='  -0.62%  '; Text=$false}
    @{R=50; C='B'; V='EnergySwap'; Text=$false}
    @{R=50; C='C'; V='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Text=$false}
    @{R=50; C='D'; V='7.53'; Text=$true}
    @{R=50; C='E'; V='  -0.64%  '; Text=$false}
    @{R=51; C='B'; V='Mantle'; Text=$false}
    @{R=51; C='C'; V='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; Text=$false}
    @{R=51; C='D'; V='0.407'; Text=$true}
    @{R=51; C='E'; V='  -0.71%  '; Text=$false}
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.R, $colIndex[$u.C])
    if ($u.Text) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.V
}
